$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Single-value cells (rows are 1-indexed in the COM model)
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "417"
$t.Cell(6, 1).Range.Text = "0.00070"
$t.Cell(7, 1).Range.Text = "0.00023"
$t.Cell(8, 1).Range.Text = "0.00007"
$t.Cell(9, 1).Range.Text = "0.00034"
$t.Cell(10, 1).Range.Text = "0.00049"
$t.Cell(11, 1).Range.Text = "0.00052"
$t.Cell(12, 1).Range.Text = "0.09590"

# Last three rows collapse their multi-run tab-separated content down to a
# single value each.
$t.Cell(44, 1).Range.Text = "99.93"
$t.Cell(45, 1).Range.Text = "0.1"
$t.Cell(46, 1).Range.Text = "128"
